# OWASP Top 10 Mapping (ES) - update 2017/2021 category names to the
# corrected/English-normalized wording, relabel "(Nuevo)" -> "(Nueva)",
# remove two duplicate/leftover connector arrows from the drawing,
# narrow column C, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (2017 list) ---------------------------------------------
$ws.Range("C7").Value  = "A01:2017-Injection"
$ws.Range("C8").Value  = "A02:2017-Broken Authentication"
$ws.Range("C9").Value  = "A03:2017-Sensitive Data Exposure"
$ws.Range("C10").Value = "A04:2017-XML External Entities (XXE)"
$ws.Range("C11").Value = "A05:2017-Broken Access Control"
$ws.Range("C12").Value = "A06:2017-Security Misconfiguration"
$ws.Range("C13").Value = "A07:2017-Cross-Site Scripting (XSS)"
$ws.Range("C14").Value = "A08:2017-Insecure Deserialization"
$ws.Range("C15").Value = "A09:2017-Using Components with Known Vulnerabilities"
$ws.Range("C16").Value = "A10:2017-Insufficient Logging & Monitoring"

# --- Column D ("(Nuevo)" -> "(Nueva)" markers) -------------------------
$ws.Range("D10").Value = "(Nueva)"
$ws.Range("D14").Value = "(Nueva)"
$ws.Range("D16").Value = "(Nueva)"

# --- Column E (2021 list) ----------------------------------------------
$ws.Range("E7").Value  = "A01:2021-Pérdida de Control de Acceso"
$ws.Range("E8").Value  = "A02:2021-Fallas Criptográficas"
$ws.Range("E9").Value  = "A03:2021-Inyección"
$ws.Range("E10").Value = "A04:2021-Diseño Inseguro"
$ws.Range("E11").Value = "A05:2021-Configuración de Seguridad Incorrecta"
$ws.Range("E12").Value = "A06:2021-Componentes Vulnerables y Desactualizados"
$ws.Range("E13").Value = "A07:2021-Fallas de Identificación y Autenticación"
$ws.Range("E14").Value = "A08:2021-Fallas en el Software y en la Integridad de los Datos"
$ws.Range("E15").Value = "A09:2021-Fallas en el Registro y Monitoreo*"
$ws.Range("E16").Value = "A10:2021-Falsificación de Solicitudes del Lado del Servidor (SSRF)*"
$ws.Range("E17").Value = "* A partir de la encuesta"

# --- Column widths -------------------------------------------------------
# Column B loses its explicit width in the diff (back to default sizing);
# column C narrows from ~54.3 to 48.5 raw width units.
$ws.Columns("B").ColumnWidth = 8.43
$ws.Columns("C").ColumnWidth = 47.666666666666664

# --- Drawing cleanup: remove the two extra/duplicate connector arrows ---
$ws.Shapes.Item("Straight Arrow Connector 11").Delete()
$ws.Shapes.Item("Straight Arrow Connector 12").Delete()

# --- Selection -----------------------------------------------------------
[void]$ws.Range("E23").Select()

Write-Host "done"
